$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 14929352
$ws.Range("I100").Value = 20834700
$ws.Range("K100").Value = 20834700
$ws.Range("M100").Value = -20834159
$ws.Range("H113").Value = 9910.883
$ws.Range("I113").Value = 2421.25
$ws.Range("J113").Value = 12215.385
$ws.Range("K113").Value = 2421.25
$ws.Range("L113").Value = 12215.385
$ws.Range("M113").Value = 832.75
$ws.Range("N113").Value = -18723.385
$ws.Range("H132").Value = 42341.32
$ws.Range("I132").Value = 43972.25
$ws.Range("K132").Value = 131916.75
$ws.Range("M132").Value = -129386.75
$ws.Range("H137").Value = 1667.6957
$ws.Range("I137").Value = 1615.3529
$ws.Range("J137").Value = 1816
$ws.Range("K137").Value = 4846.0587
$ws.Range("L137").Value = 5448
$ws.Range("M137").Value = -2296.0587
$ws.Range("N137").Value = -10548
$ws.Range("H138").Value = 12051058
$ws.Range("I138").Value = 1308.1951
$ws.Range("J138").Value = 23813908
$ws.Range("K138").Value = 3924.5853
$ws.Range("L138").Value = 71441724
$ws.Range("M138").Value = 1215.4147
$ws.Range("N138").Value = -71452004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3520.4824
$ws.Range("I32").Value = 2319.2896
$ws.Range("K32").Value = 2319.2896
$ws.Range("M32").Value = -2032.2896
$ws.Range("H88").Value = 6638.0835
$ws.Range("I88").Value = 10405.667
$ws.Range("J88").Value = 2870.5
$ws.Range("K88").Value = 10405.667
$ws.Range("L88").Value = 2870.5
$ws.Range("M88").Value = -9999.666999999999
$ws.Range("N88").Value = -3682.5
$ws.Range("H91").Value = 6638.0835
$ws.Range("I91").Value = 10405.667
$ws.Range("J91").Value = 2870.5
$ws.Range("K91").Value = 10405.667
$ws.Range("L91").Value = 2870.5
$ws.Range("M91").Value = -9001.666999999999
$ws.Range("N91").Value = -5678.5
$ws.Range("H97").Value = 1331.4615
$ws.Range("I97").Value = 1162.9
$ws.Range("J97").Value = 1893.3334
$ws.Range("K97").Value = 1162.9
$ws.Range("L97").Value = 1893.3334
$ws.Range("M97").Value = -666.9000000000001
$ws.Range("N97").Value = -2885.3334
$ws.Range("H107").Value = 44381.668
$ws.Range("J107").Value = 44381.668
$ws.Range("L107").Value = 44381.668
$ws.Range("N107").Value = -52061.668
$ws.Range("H109").Value = 18825.834
$ws.Range("J109").Value = 18825.834
$ws.Range("L109").Value = 18825.834
$ws.Range("N109").Value = -21599.834

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1971.7059
$ws.Range("I86").Value = 2001.2
$ws.Range("J86").Value = 1929.5714
$ws.Range("K86").Value = 2001.2
$ws.Range("L86").Value = 1929.5714
$ws.Range("M86").Value = -878.2
$ws.Range("N86").Value = -4175.5714
$ws.Range("H89").Value = 1971.7059
$ws.Range("I89").Value = 2001.2
$ws.Range("J89").Value = 1929.5714
$ws.Range("K89").Value = 10006
$ws.Range("L89").Value = 9647.857
$ws.Range("M89").Value = -4390
$ws.Range("N89").Value = -20879.857
$ws.Range("H94").Value = 547.3889
$ws.Range("I94").Value = 531.94116
$ws.Range("K94").Value = 531.94116
$ws.Range("M94").Value = -80.94115999999997
$ws.Range("H99").Value = 1471
$ws.Range("I99").Value = 883.2222
$ws.Range("K99").Value = 883.2222
$ws.Range("M99").Value = 614.7778
$ws.Range("H107").Value = 961.625
$ws.Range("I107").Value = 970.4286
$ws.Range("K107").Value = 970.4286
$ws.Range("M107").Value = 949.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 32295
$ws.Range("J68").Value = 32295
$ws.Range("L68").Value = 32295
$ws.Range("N68").Value = -33793
$ws.Range("H71").Value = 32295
$ws.Range("J71").Value = 32295
$ws.Range("L71").Value = 96885
$ws.Range("N71").Value = -104373
$ws.Range("H74").Value = 18588
$ws.Range("J74").Value = 18588
$ws.Range("L74").Value = 18588
$ws.Range("N74").Value = -20336
$ws.Range("H77").Value = 18588
$ws.Range("J77").Value = 18588
$ws.Range("L77").Value = 55764
$ws.Range("N77").Value = -64500
$ws.Range("H100").Value = 37820
$ws.Range("J100").Value = 37820
$ws.Range("L100").Value = 37820
$ws.Range("N100").Value = -39984
$ws.Range("H123").Value = 47978.332
$ws.Range("J123").Value = 47978.332
$ws.Range("L123").Value = 47978.332
$ws.Range("N123").Value = -57778.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 670742.6
$ws.Range("I68").Value = 1201042.4
$ws.Range("J68").Value = 1554.8572
$ws.Range("K68").Value = 3603127.2
$ws.Range("L68").Value = 4664.571599999999
$ws.Range("M68").Value = -3602316.2
$ws.Range("N68").Value = -6286.571599999999
$ws.Range("H71").Value = 670742.6
$ws.Range("I71").Value = 1201042.4
$ws.Range("J71").Value = 1554.8572
$ws.Range("K71").Value = 10809381.6
$ws.Range("L71").Value = 13993.7148
$ws.Range("M71").Value = -10805325.6
$ws.Range("N71").Value = -22105.7148
$ws.Range("H107").Value = 767.44183
$ws.Range("I107").Value = 678.53845
$ws.Range("J107").Value = 903.41174
$ws.Range("K107").Value = 2035.61535
$ws.Range("L107").Value = 2710.23522
$ws.Range("M107").Value = -115.61535
$ws.Range("N107").Value = -6550.23522
$ws.Range("H113").Value = 1183.6957
$ws.Range("I113").Value = 1183.08
$ws.Range("K113").Value = 3549.24
$ws.Range("M113").Value = -1379.24
$ws.Range("H119").Value = 4950.0625
$ws.Range("I119").Value = 2361.2856
$ws.Range("J119").Value = 6963.5557
$ws.Range("K119").Value = 7083.8568
$ws.Range("L119").Value = 20890.6671
$ws.Range("M119").Value = -2245.8568
$ws.Range("N119").Value = -30566.6671
$ws.Range("H131").Value = 3541.625
$ws.Range("I131").Value = 622.1111
$ws.Range("J131").Value = 4215.359
$ws.Range("K131").Value = 1866.3333
$ws.Range("L131").Value = 12646.077
$ws.Range("M131").Value = 3173.6667
$ws.Range("N131").Value = -22726.077
$ws.Range("H132").Value = 1421.5938
$ws.Range("I132").Value = 1608.875
$ws.Range("J132").Value = 1234.3125
$ws.Range("K132").Value = 14479.875
$ws.Range("L132").Value = 11108.8125
$ws.Range("M132").Value = -11949.875
$ws.Range("N132").Value = -16168.8125
$ws.Range("H137").Value = 2370.9412
$ws.Range("J137").Value = 3097.5334
$ws.Range("L137").Value = 9292.600199999999
$ws.Range("N137").Value = -19492.6002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 40000
$ws.Range("J51").Value = 40000
$ws.Range("L51").Value = 40000
$ws.Range("N51").Value = -41018

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 939.6667
$ws.Range("I93").Value = 963
$ws.Range("J93").Value = 893
$ws.Range("K93").Value = 963
$ws.Range("L93").Value = 893
$ws.Range("M93").Value = 285
$ws.Range("N93").Value = -3389
